# B6-PowerPoint.pptx edit
#
# 1) Three tables (on slides 14, 15, 16) switch from the deck's single
#    custom "Table_0" style ({7F1115E1-6101-4312-B7F5-166261ABD05F}) to
#    PowerPoint's built-in "No Style, No Grid" table style
#    ({05152B13-D281-4F62-BA44-350C975AB2C9}).
# 2) The presentation's theme is changed from "Integral" (Red Violet
#    colour scheme) to the built-in "Office Theme" colour scheme.

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables ----------------------------------------
$newStyleId = "{05152B13-D281-4F62-BA44-350C975AB2C9}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Re-colour the theme (Integral / Red Violet -> Office Theme) ------
# Theme colour scheme slots, in PowerPoint's fixed order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# RGB is packed as r + g*256 + b*65536 (PowerPoint's native ordering).
$officeThemeRgb = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeRgb[$i - 1]
}
